$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Update training-plan rows (columns B, C, E) with new course content ---
$ws2.Range("B4").Value = "公司业务机器中表现层的交互逻辑，UI设计"
$ws2.Range("C4").Value = "让学生在了解公司业务流程的前提下，应用前端UI交互的知识，依据要求对业务中的需求进行实现。"
$ws2.Range("E4").Value = "让该生通过对Web三板斧（HTML、CSS、JavaScript）的培训学习，并在其基础上，利用前端MVC框架，MVVM框架对具体的业务进行实现，产出价值。"

$ws2.Range("B5").Value = "职场沟通方式及沟通技巧"
$ws2.Range("C5").Value = "让学生了解职场的沟通方式（电子邮件、电话沟通、PPT演讲等），训练学生的沟通礼仪和技巧，培养表达和演讲能力。"
$ws2.Range("E5").Value = "定期让该生进行PPT演讲，锻炼了其表达和演讲能力。在工作中，该生也能积极地与同事进行沟通交流，"

$ws2.Range("B6").Value = "职场办公软件的应用"
$ws2.Range("C6").Value = "让学生熟悉MircoSoft Offic办公套件的使用，适应公司制式化业务办公中的任务交接及安排、汇报等流程"
$ws2.Range("E6").Value = "在对应的技术培训期间要求该生总结工作以及学习报告并交纳对应的负责人，去糟粕取精华，并与公司现存的经验碰撞催化产出。"

$ws2.Range("B7").Value = "项目版本管理和版本迭代"
$ws2.Range("C7").Value = "让学生在实际的项目开发过程中逐渐了解并行开发的流程，以及基于Git、SVN等工具的版本管理和Code Review、Pull Request的流程"
$ws2.Range("E7").Value = "通过在项目中的开发过程中负责人的follow下，逐步学习Git等版本工具的使用，理解Git常用命令的作用和应用的场景，项目版本错乱时的解决手段等等。"

$ws2.Range("B8").Value = "底层服务器部署以及Linux系统的学习"
$ws2.Range("C8").Value = "主要了解和熟悉基于Linux、Windows平台的环境部署、安全配置等，了解底层服务器中的REST API以及Debug。"
$ws2.Range("E8").Value = "该生通过手动部署泛微OA系统、筑店商品管理系统熟悉了Windows和Linux下的环境搭建，也了解了NGINX和APACHE的相关配置，能够熟悉使用宝塔面板对Linux进行基本运维。"

$ws2.Range("B9").Value = "Node环境下Modern JavaScript生态"
$ws2.Range("E9").Value = "通过在公司现有“自动化”工具项目中通过Vue生态下的相关技术知识、配合培训的要求完成前端部分相关的设计、开发工作"
$ws2.Range("C9").Value = "包括但不限于Webpack等前端部署工具、Grunt、Gulp、Babel等构建工具、TypeScript下的OOP风格的设计模式、Jshint/Prettier等辅助开发工具的学习和练习"

# --- Clear the now-unused duplicate G/H columns for the header and data rows (keep formatting) ---
$ws2.Range("G3:H9").ClearContents()

# --- Remove the last two (now superfluous) data rows entirely ---
$ws2.Rows("10:11").Delete()

# --- Remove the trailing duplicate G/H cells (incl. formulas) on the total row ---
$ws2.Range("G10:H10").Clear()

# --- Restore sheet2 view state (frozen pane / selection) to match the edited layout ---
$ws2.Application.ActiveWindow.ScrollColumn = 1
$ws2.Range("F11").Select()
